$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row (previously only "Status" existed in A1; now split across A1:C1)
$ws.Range("A1").Value = "Sr. No."
$ws.Range("B1").Value = "Test Case Name"
$ws.Range("C1").Value = "Status"

# Row 2: testFileUpload, now PASS
$ws.Range("B2").Value = "testFileUpload"
$ws.Range("C2").Value = "PASS"

# Row 3: testABTestingLinkNavigation, now PASS
$ws.Range("B3").Value = "testABTestingLinkNavigation"
$ws.Range("C3").Value = "PASS"

# Row 4: testLinksCountOnHomePage, now PASS
$ws.Range("B4").Value = "testLinksCountOnHomePage"
$ws.Range("C4").Value = "PASS"
